$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 3.15
$ws.Range("L2").Value = 3.7
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.98
$ws.Range("U2").Value = 1.6
$ws.Range("V2").Value = 2.22
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 22
$ws.Range("AD2").Value = 6.3
$ws.Range("AE2").Value = 11.75
$ws.Range("AH2").Value = 10.25
$ws.Range("AI2").Value = 17
$ws.Range("AK2").Value = 40
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 30
$ws.Range("AP2").Value = 17
$ws.Range("AU2").Value = 6.5
$ws.Range("AW2").Value = 5.2
$ws.Range("AX2").Value = 17.5
